$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Apple Watch"
$ws.Range("B3").Value = "Electronics"

$ws.Columns("A:A").AutoFit() | Out-Null

$ws.Range("B3").Select()
